# Updated symbol list on Tue Dec 20 19:26:04 UTC 2022 with GitHub Actions
# Applies updated Price (column D) and Volume(1h) (column E) values
# for the crypto symbol list snapshot.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell -> new value. Numeric-looking "Price" values are written
# with a leading apostrophe so they stay text cells (matching the
# original inline-string/text storage) instead of being coerced to
# numbers.
$updates = [ordered]@{
    "D2"  = "'249.96"
    "D3"  = "'22.91"
    "D4"  = "'5.421"
    "D5"  = "'0.05647"
    "D6"  = "'3.422"
    "D7"  = "'6.367"
    "D8"  = "'0.8147"
    "D9"  = "'0.9138"
    "D10" = "'0.0005788"
    "E10" = "9OneONE"
    "D11" = "'0.1441"
    "D12" = "'0.07533"
    "D13" = "'0.03126"
    "D14" = "'0.03109"
    "D15" = "'0.09349"
    "D16" = "'3.561"
    "D17" = "'0.001588"
    "D18" = "'0.04756"
    "D20" = "'0.004999"
    "D21" = "'0.001032"
    "D28" = "'0.0003030"
    "D40" = "'0.04043"
    "D41" = "'0.006838"
    "D42" = "'0.1069"
    "D44" = "'0.007550"
    "D48" = "'0.2412"
    "E48" = "47BOLOBOLOBestin24h"
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
